$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AStar")

$data = @(
    @('A*', 0.0977, 84, 57, 'LOPEZ', 'AARON, LEGAL, PEERS, BLEND, CHEAP, CLERK, DELAY, EQUAL, IDEAL, MEDAL, METAL, SPEAK, SPEND, SPENT, SPERM, SPINE, STEAL, EVENT, EVERY, SCENE, TEENS, THERE, WHERE, LUCAS, LYING, BEING, BLIND, BLINK, CHEAT, ESSAY, FIBRE, IDEAS, ISLAM, MYERS, QUERY, SHINE, STEAM, TEXAS, USERS, USUAL, VEGAS, VITAL, WEIRD, WHEAT, DUBAI, DYING, SHIRT, SIGNS, SKINS, SKIRT, SQUAD, SWING, THING, THINK, THIRD, TWINS, USING'),
    @('A*', 0.111, 84, 57, 'LOPEZ', 'AARON, LEGAL, PEERS, BLEND, CHEAP, CLERK, DELAY, EQUAL, IDEAL, MEDAL, METAL, SPEAK, SPEND, SPENT, SPERM, SPINE, STEAL, EVENT, EVERY, SCENE, TEENS, THERE, WHERE, LUCAS, LYING, BEING, BLIND, BLINK, CHEAT, ESSAY, FIBRE, IDEAS, ISLAM, MYERS, QUERY, SHINE, STEAM, TEXAS, USERS, USUAL, VEGAS, VITAL, WEIRD, WHEAT, DUBAI, DYING, SHIRT, SIGNS, SKINS, SKIRT, SQUAD, SWING, THING, THINK, THIRD, TWINS, USING'),
    @('A*', 0.09719999999999999, 84, 57, 'LOPEZ', 'AARON, LEGAL, PEERS, BLEND, CHEAP, CLERK, DELAY, EQUAL, IDEAL, MEDAL, METAL, SPEAK, SPEND, SPENT, SPERM, SPINE, STEAL, EVENT, EVERY, SCENE, TEENS, THERE, WHERE, LUCAS, LYING, BEING, BLIND, BLINK, CHEAT, ESSAY, FIBRE, IDEAS, ISLAM, MYERS, QUERY, SHINE, STEAM, TEXAS, USERS, USUAL, VEGAS, VITAL, WEIRD, WHEAT, DUBAI, DYING, SHIRT, SIGNS, SKINS, SKIRT, SQUAD, SWING, THING, THINK, THIRD, TWINS, USING'),
    @('A*', 0.003, 233, 2, 'LIKED', 'AARON, LIKED'),
    @('A*', 0.1069, 101, 57, 'POSTS', 'AARON, SPENT, PEERS, ESSAY, SIGNS, SKINS, TEENS, TEXAS, TWINS, USERS, SHIRT, SKIRT, SPEAK, SPEND, SPERM, SPINE, STEAL, STEAM, IDEAS, LUCAS, MYERS, VEGAS, CHEAP, CHEAT, EVENT, ISLAM, METAL, SCENE, SHINE, SQUAD, SWING, THERE, THING, THINK, THIRD, USING, USUAL, VITAL, WHEAT, BEING, BLEND, BLIND, BLINK, CLERK, DELAY, DUBAI, DYING, EQUAL, EVERY, FIBRE, IDEAL, LEGAL, LYING, MEDAL, QUERY, WEIRD, WHERE'),
    @('A*', 0.0322, 19, 19, 'DAIRY', 'AARON, DANNY, CANDY, HANDY, SANDY, DANCE, FANCY, BANDS, HANDS, LANDS, PANIC, BANKS, JANET, LANCE, LANES, PANEL, PANTS, TANKS, WANTS'),
    @('A*', 0.002, 1, 1, 'APART', 'AARON'),
    @('A*', 0.0135, 11, 9, 'ALIAS', 'AARON, ARISE, ANIME, ANGEL, ANGLE, ANNIE, ANNEX, ARGUE, ARMED'),
    @('A*', 0.0186, 12, 12, 'INTER', 'AARON, TUNER, BLAIR, CHAIR, CHOIR, BINGO, KENYA, KINDA, LINDA, FLOUR, IDAHO, PLAZA'),
    @('A*', 0.09370000000000001, 107, 57, 'COMIC', 'AARON, ISLAM, CHEAP, CHEAT, CLERK, BEING, BLIND, BLINK, DUBAI, DYING, FIBRE, IDEAL, IDEAS, LUCAS, LYING, MEDAL, METAL, MYERS, SCENE, SHINE, SHIRT, SIGNS, SKINS, SKIRT, SPERM, SPINE, STEAM, SWING, THING, THINK, THIRD, TWINS, USING, VITAL, WEIRD, BLEND, DELAY, EQUAL, ESSAY, EVENT, EVERY, LEGAL, PEERS, QUERY, SPEAK, SPEND, SPENT, SQUAD, STEAL, TEENS, TEXAS, THERE, USERS, USUAL, VEGAS, WHEAT, WHERE'),
    @('A*', 0.003, 233, 2, 'SILLY', 'AARON, SILLY'),
    @('A*', 0.002, 233, 2, 'FEELS', 'AARON, FEELS'),
    @('A*', 0.0335, 22, 21, 'LEAVE', 'AARON, OLIVE, REBEL, NEEDS, NEWLY, REFER, REPLY, RESET, RIVER, RULED, RULES, OXIDE, REMIX, RIDGE, RIDER, RIDES, NIGHT, OUGHT, RIGHT, RISKS, RUGBY'),
    @('A*', 0.032, 22, 21, 'LEAVE', 'AARON, OLIVE, REBEL, NEEDS, NEWLY, REFER, REPLY, RESET, RIVER, RULED, RULES, OXIDE, REMIX, RIDGE, RIDER, RIDES, NIGHT, OUGHT, RIGHT, RISKS, RUGBY'),
    @('A*', 0.1002, 71, 62, 'DROPS', 'AARON, SPORT, SWORD, SCORE, SHORE, SHORT, STORE, STORM, STORY, SHARP, SPARC, SPARE, STARS, ICONS, LIONS, TIONS, BEARS, BLOND, DIARY, FEARS, GLORY, IVORY, PHONE, PLANS, STONE, TEARS, WHORE, YEARS, GUARD, HEARD, PEARL, SCARY, SHARE, SHARK, SMART, SPANK, STAND, START, BEANS, CLONE, DIANE, EBONY, EVANS, FLOAT, JEANS, LEONE, MEANS, THONG, CHARM, CHART, CLARK, GENRE, HEART, HENRY, ISAAC, PLANE, PLANT, BLANK, FINAL, GIANT, MEANT, THANK')
)

$startRow = 12
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}
